$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra "Hora de estudar" automation row entirely.
$ws.Rows(4).Delete()

# Drop the stray empty F1 cell left behind by the dynamic button automation.
$ws.Range("F1").Clear()

# Fix the values the buggy automation had been writing wrong.
$ws.Range("A3").Value = "Boa tarde"
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = $false
$ws.Range("F3").Value = $false

# The old fixed row heights (19.5 / 18.75) were a leftover of the buggy
# dynamic button automation; let the rows size themselves again.
$ws.Rows("1:3").AutoFit()
